# Append timed-response rows (134-167) captured from the uploaded text-file
# questions, matching the existing "Question / Model Name / Response" layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(134, 1).Value2 = "What's the maximum number of lithology types in an log?"
$ws.Cells.Item(134, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(134, 3).Value2 = "The highest number of lithology types that can be represented in a log is 450."

$ws.Cells.Item(135, 1).Value2 = "How many tracks can you define in one ODF?"
$ws.Cells.Item(135, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(135, 3).Value2 = "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."

$ws.Cells.Item(136, 1).Value2 = "How many curve shades can I create?"
$ws.Cells.Item(136, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(136, 3).Value2 = "According to the document, you can create 250 curve shades."

$ws.Cells.Item(137, 1).Value2 = "How many curves can I load in one go?"
$ws.Cells.Item(137, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(137, 3).Value2 = "You can load up to 450 curves at a time."

$ws.Cells.Item(138, 1).Value2 = "What the maximum number of headers I can display in my log?"
$ws.Cells.Item(138, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(138, 3).Value2 = "The maximum number of headers you can display in your log is up to 50."

$ws.Cells.Item(139, 1).Value2 = "How many tables can I have in my log?"
$ws.Cells.Item(139, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(139, 3).Value2 = "According to the GEO application documentation, there is no specific limit on the number of tables that can be used in a composite log. However, it's worth noting that some systems, such as the corporate geological database, may require certain information to be cataloged separately.
If you need more information or clarification, please refer to the `"Tables`" section under the title `"---Feedback---`"."

$ws.Cells.Item(140, 1).Value2 = "Whats the maximum number of characters in a single text entry?"
$ws.Cells.Item(140, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(140, 3).Value2 = "The maximum number of characters for an individual cell is 1999."

$ws.Cells.Item(141, 1).Value2 = "How many symbols can I have in the plot at any one time?"
$ws.Cells.Item(141, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(141, 3).Value2 = "You can have up to 10,000 symbols in a plot at any given time."

$ws.Cells.Item(142, 1).Value2 = "How many scales can I define?"
$ws.Cells.Item(142, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(142, 3).Value2 = "According to the document, you can define up to 23 scales."

$ws.Cells.Item(143, 1).Value2 = "What the maximum number of data files I can load?"
$ws.Cells.Item(143, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(143, 3).Value2 = "The maximum number of data files you can load is unlimited."

$ws.Cells.Item(144, 1).Value2 = "What's the maximum number of lithology types in an log?"
$ws.Cells.Item(144, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(144, 3).Value2 = "The highest number of lithology types that can be represented in a log is 450."

$ws.Cells.Item(145, 1).Value2 = "How many tracks can you define in one ODF?"
$ws.Cells.Item(145, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(145, 3).Value2 = "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."

$ws.Cells.Item(146, 1).Value2 = "How many curve shades can I create?"
$ws.Cells.Item(146, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(146, 3).Value2 = "According to the document, you can create 250 curve shades."

$ws.Cells.Item(147, 1).Value2 = "How many curves can I load in one go?"
$ws.Cells.Item(147, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(147, 3).Value2 = "You can load up to 450 curves at a time."

$ws.Cells.Item(148, 1).Value2 = "What's the maximum number of lithology types in an log?"
$ws.Cells.Item(148, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(148, 3).Value2 = "The highest number of lithology types that can be represented in a log is 450."

$ws.Cells.Item(149, 1).Value2 = "How many tracks can you define in one ODF?"
$ws.Cells.Item(149, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(149, 3).Value2 = "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."

$ws.Cells.Item(150, 1).Value2 = "How many curve shades can I create?"
$ws.Cells.Item(150, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(150, 3).Value2 = "According to the document, you can create 250 curve shades."

$ws.Cells.Item(151, 1).Value2 = "How many curves can I load in one go?"
$ws.Cells.Item(151, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(151, 3).Value2 = "You can load up to 450 curves at a time."

$ws.Cells.Item(152, 1).Value2 = "What the maximum number of headers I can display in my log?"
$ws.Cells.Item(152, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(152, 3).Value2 = "The maximum number of headers you can display in your log is up to 50."

$ws.Cells.Item(153, 1).Value2 = "How many tables can I have in my log?"
$ws.Cells.Item(153, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(153, 3).Value2 = "According to the GEO application documentation, there is no specific limit on the number of tables that can be used in a composite log. However, it's worth noting that some systems, such as the corporate geological database, may require certain information to be cataloged separately.
If you need more information or clarification, please refer to the `"Tables`" section under the title `"---Feedback---`"."

$ws.Cells.Item(154, 1).Value2 = "Whats the maximum number of characters in a single text entry?"
$ws.Cells.Item(154, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(154, 3).Value2 = "The maximum number of characters for an individual cell is 1999."

$ws.Cells.Item(155, 1).Value2 = "How many symbols can I have in the plot at any one time?"
$ws.Cells.Item(155, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(155, 3).Value2 = "You can have up to 10,000 symbols in a plot at any given time."

$ws.Cells.Item(156, 1).Value2 = "How many scales can I define?"
$ws.Cells.Item(156, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(156, 3).Value2 = "According to the document, you can define up to 23 scales."

$ws.Cells.Item(157, 1).Value2 = "What the maximum number of data files I can load?"
$ws.Cells.Item(157, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(157, 3).Value2 = "The maximum number of data files you can load is unlimited."

$ws.Cells.Item(158, 1).Value2 = "What's the maximum number of lithology types in an log?"
$ws.Cells.Item(158, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(158, 3).Value2 = "The highest number of lithology types that can be represented in a log is 450."

$ws.Cells.Item(159, 1).Value2 = "How many tracks can you define in one ODF?"
$ws.Cells.Item(159, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(159, 3).Value2 = "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."

$ws.Cells.Item(160, 1).Value2 = "How many curve shades can I create?"
$ws.Cells.Item(160, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(160, 3).Value2 = "According to the document, you can create 250 curve shades."

$ws.Cells.Item(161, 1).Value2 = "How many curves can I load in one go?"
$ws.Cells.Item(161, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(161, 3).Value2 = "You can load up to 450 curves at a time."

$ws.Cells.Item(162, 1).Value2 = "What the maximum number of headers I can display in my log?"
$ws.Cells.Item(162, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(162, 3).Value2 = "The maximum number of headers you can display in your log is up to 50."

$ws.Cells.Item(163, 1).Value2 = "How many tables can I have in my log?"
$ws.Cells.Item(163, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(163, 3).Value2 = "According to the GEO application documentation, there is no specific limit on the number of tables that can be used in a composite log. However, it's worth noting that some systems, such as the corporate geological database, may require certain information to be cataloged separately.
If you need more information or clarification, please refer to the `"Tables`" section under the title `"---Feedback---`"."

$ws.Cells.Item(164, 1).Value2 = "Whats the maximum number of characters in a single text entry?"
$ws.Cells.Item(164, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(164, 3).Value2 = "The maximum number of characters for an individual cell is 1999."

$ws.Cells.Item(165, 1).Value2 = "How many symbols can I have in the plot at any one time?"
$ws.Cells.Item(165, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(165, 3).Value2 = "You can have up to 10,000 symbols in a plot at any given time."

$ws.Cells.Item(166, 1).Value2 = "How many scales can I define?"
$ws.Cells.Item(166, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(166, 3).Value2 = "According to the document, you can define up to 23 scales."

$ws.Cells.Item(167, 1).Value2 = "What the maximum number of data files I can load?"
$ws.Cells.Item(167, 2).Value2 = "llama3.2:latest"
$ws.Cells.Item(167, 3).Value2 = "The maximum number of data files you can load is unlimited."
